# Apply updated Betfair back/lay odds values to Sheet1.
# Values below reflect the refreshed odds snapshot for 2026-01-11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.75
$ws.Range("G2").Value = 1.88
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 4.2
$ws.Range("K2").Value = 4.4
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 4.4
$ws.Range("P2").Value = 2.28
$ws.Range("Q2").Value = 1.56
$ws.Range("S2").Value = 2.38
$ws.Range("T2").Value = 1.01
$ws.Range("V2").Value = 1.25
$ws.Range("W2").Value = 2.14
$ws.Range("X2").Value = 21
$ws.Range("Y2").Value = 24
$ws.Range("Z2").Value = 42
$ws.Range("AA2").Value = 110
$ws.Range("AB2").Value = 12.5
$ws.Range("AC2").Value = 9.800000000000001
$ws.Range("AD2").Value = 20
$ws.Range("AE2").Value = 55
$ws.Range("AF2").Value = 15
$ws.Range("AG2").Value = 10.5
$ws.Range("AH2").Value = 17.5
$ws.Range("AI2").Value = 55
$ws.Range("AJ2").Value = 21
$ws.Range("AK2").Value = 17.5
$ws.Range("AL2").Value = 28
$ws.Range("AM2").Value = 80
$ws.Range("AN2").Value = 9.199999999999999
$ws.Range("AO2").Value = 46
# Row 3
$ws.Range("F3").Value = 1.72
$ws.Range("L3").Value = 1.25
$ws.Range("R3").Value = 1.73
$ws.Range("T3").Value = 1.55
# Row 5
$ws.Range("F5").Value = 1.2
$ws.Range("G5").Value = 1.8
$ws.Range("H5").Value = 2.24
$ws.Range("J5").Value = 2.24
# Row 8
$ws.Range("F8").Value = 1.88
$ws.Range("G8").Value = 2.1
$ws.Range("H8").Value = 3.6
$ws.Range("P8").Value = 2.4
# Row 9
$ws.Range("F9").Value = 1.09
$ws.Range("G9").Value = 1000
$ws.Range("H9").Value = 1.04
$ws.Range("I9").Value = 1000
$ws.Range("J9").Value = 1.01
$ws.Range("K9").Value = 980
$ws.Range("P9").Value = 1.25
$ws.Range("Q9").Value = 1.01
# Row 10
$ws.Range("F10").Value = 3
$ws.Range("J10").Value = 3.9
$ws.Range("P10").Value = 2.52
$ws.Range("Q10").Value = 1.46
# Row 14
$ws.Range("F14").Value = 2.8
$ws.Range("H14").Value = 3.1
$ws.Range("AN14").Value = 140
# Row 15
$ws.Range("P15").Value = 2.08
# Row 19
$ws.Range("G19").Value = 2.76
$ws.Range("I19").Value = 3.45
$ws.Range("J19").Value = 3
# Row 25
$ws.Range("Q25").Value = 1.69
# Row 27
$ws.Range("F27").Value = 1.24
$ws.Range("G27").Value = 1000
$ws.Range("H27").Value = 1.53
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1.17
$ws.Range("P27").Value = 1.25
$ws.Range("Q27").Value = 1.02
# Row 30
$ws.Range("F30").Value = 3.95
$ws.Range("G30").Value = 4.1
$ws.Range("K30").Value = 3.65
# Row 31
$ws.Range("F31").Value = 3.3
$ws.Range("G31").Value = 3.75
$ws.Range("I31").Value = 2.34
$ws.Range("J31").Value = 3.5
$ws.Range("K31").Value = 3.65
$ws.Range("P31").Value = 2
$ws.Range("Q31").Value = 1.83
# Row 32
$ws.Range("I32").Value = 1.53
$ws.Range("P32").Value = 1.87
# Row 33
$ws.Range("Q33").Value = 1.91
$ws.Range("R33").Value = 1.41
$ws.Range("AO33").Value = 55
# Row 35
$ws.Range("H35").Value = 1.22
$ws.Range("K35").Value = 7.4
# Row 37
$ws.Range("F37").Value = 2.92
$ws.Range("G37").Value = 3.25
$ws.Range("H37").Value = 2.42
$ws.Range("I37").Value = 2.7
$ws.Range("J37").Value = 3.3
$ws.Range("K37").Value = 3.7
$ws.Range("P37").Value = 1.89
$ws.Range("Q37").Value = 1.94
# Row 38
$ws.Range("F38").Value = 2.1
$ws.Range("H38").Value = 3.75
$ws.Range("I38").Value = 4.3
# Row 40
$ws.Range("G40").Value = 1.95
$ws.Range("Q40").Value = 1.5
# Row 41
$ws.Range("I41").Value = 24
$ws.Range("J41").Value = 11.5
$ws.Range("K41").Value = 12.5
$ws.Range("Q41").Value = 1.2
$ws.Range("R41").Value = 2.72
$ws.Range("S41").Value = 1.54
$ws.Range("U41").Value = 2.1
# Row 47
$ws.Range("P47").Value = 2.08
$ws.Range("Q47").Value = 1.76
$ws.Range("AD47").Value = 10.5
$ws.Range("AK47").Value = 75
# Row 48
$ws.Range("Q48").Value = 2.66
$ws.Range("Z48").Value = 15
# Row 49
$ws.Range("N49").Value = 3.05
# Row 50
$ws.Range("F50").Value = 2.74
$ws.Range("G50").Value = 3.2
$ws.Range("H50").Value = 2.5
$ws.Range("I50").Value = 2.86
$ws.Range("J50").Value = 3.25
$ws.Range("K50").Value = 3.8
$ws.Range("P50").Value = 1.84
$ws.Range("Q50").Value = 1.94
# Row 51
$ws.Range("K51").Value = 5.1
# Row 52
$ws.Range("L52").Value = 1.2
$ws.Range("N52").Value = 5.9
$ws.Range("P52").Value = 2.68
$ws.Range("R52").Value = 1.67
$ws.Range("U52").Value = 1.45
$ws.Range("X52").Value = 980
$ws.Range("AC52").Value = 980
# Row 54
$ws.Range("AJ54").Value = 18.5

$wb.Save()
